$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.892.85'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '2.054.00'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.33'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.615'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.20'
$ws.Range('E8').Value = '  +3.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0801'
$ws.Range('E10').Value = '  +1.71%  '
$ws.Range('E11').Value = '  -0.80%  '
$ws.Range('D12').Value = '2.362.18'
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.48'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.59'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.25'
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.746'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '2.054.44'
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('D18').Value = '37.847.91'
$ws.Range('E18').Value = '  +2.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.21'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.75'
$ws.Range('E20').Value = '  +1.11%  '
$ws.Range('D21').Value = '0.0₃0828'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.69'
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.45'
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('E25').Value = '  +2.89%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.83'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('E28').Value = '  +6.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.06'
$ws.Range('E29').Value = '  +1.45%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.53'
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.54'
$ws.Range('E34').Value = '  +2.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.00'
$ws.Range('E35').Value = '  +8.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.35'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.97'
$ws.Range('E37').Value = '  +10.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.28'
$ws.Range('E38').Value = '  +3.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '99.05'
$ws.Range('E40').Value = '  +3.71%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0218'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.480.73'
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0950'
$ws.Range('E43').Value = '  +2.61%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.85'
$ws.Range('E44').Value = '  +2.12%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.63'
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.12'
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('E47').Value = '  +15.64%  '
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.15'
$ws.Range('E49').Value = '  -2.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.94'
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').Value = '2.249.75'
$ws.Range('E51').Value = '  +1.95%  '
